$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10 previously held 18; restore/update it to 1 per the target revision.
$ws.Range("C10").Value = 1
